# Add an "Italy" test-data sheet, modeled on the existing "Germany" sheet,
# and make it the active tab (mirrors the other per-country worksheets in
# this workbook: same layout/styles, market name + part number updated).

$wb = $excel.ActiveWorkbook

# --- Germany sheet: selection becomes "select all" (whole sheet) ---------
$germany = $wb.Worksheets.Item("Germany")
$germany.Cells.Select() | Out-Null

# --- Create "Italy" as a copy of "Germany", placed after "Slovakia" ------
$slovakia = $wb.Worksheets.Item("Slovakia")
$germany.Copy($null, $slovakia)
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Germany has an extra "FAT-S" accessories row (row 12) that the other
# country sheets (and Italy) don't carry - remove it so the layout matches
# the rest of the workbook (13 data rows instead of 14).
$italy.Rows.Item(12).Delete() | Out-Null

# --- Country-specific values ----------------------------------------------
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2454/T2453"

# --- Make Italy the active sheet/tab, cursor left on B14 -----------------
$italy.Activate() | Out-Null
$italy.Range("B14").Select() | Out-Null
